# Apply metadata updates to the "Metadata" worksheet of the
# MindfulnessSettingCS CodeSystem workbook:
#   1. Refresh the generation "Date" timestamp.
#   2. Fill in the previously-empty "Description" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 8 -> Property "Date": update the generated timestamp value.
$ws.Range("B8").Value = "2025-10-02T18:31:12+01:00"

# Row 17 -> Property "Description": set the (previously blank) description text.
$ws.Range("B17").Value = "CodeSystem defining different settings where mindfulness practice can occur"
